# Importing the Processed data files into python
#
# The "Stations" sheet had two blank placeholder rows (rows 2 and 3) above
# the data block. Remove them so the data (previously starting at row 4)
# shifts up to start at row 2, shrinking the used range from A1:F106 to
# A1:F104.

$wb = $excel.ActiveWorkbook
$stations = $wb.Worksheets.Item("Stations")
$trainServices = $wb.Worksheets.Item("TrainServices")

# Delete the two empty rows right under the header.
$stations.Rows("2:3").Delete()

# The edited file was last saved with the Stations tab active/selected,
# scrolled back to the top, with J17 selected.
$stations.Activate()
$stations.Range("J17").Select()
